# The workbook tracks daily Hortaliza (Acelga) prices. This edit inserts a
# new weekly sampling pair of rows (Primera / Segunda) at row 116, pushing
# all subsequent rows down by two positions (116->118, 117->119, ... 247->249).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 116 (each Insert() call pushes down by one row).
$ws.Rows.Item(116).Insert()
$ws.Rows.Item(116).Insert()

# New row 116: "Primera" quality sample
$ws.Cells.Item(116, 1).Value = 8
$ws.Cells.Item(116, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(116, 3).Value = "Coquimbo"
$ws.Cells.Item(116, 4).Value = 44494
$ws.Cells.Item(116, 5).Value = 4
$ws.Cells.Item(116, 6).Value = 100112009
$ws.Cells.Item(116, 7).Value = "Acelga"
$ws.Cells.Item(116, 8).Value = "Sin especificar"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 2000
$ws.Cells.Item(116, 11).Value = 500
$ws.Cells.Item(116, 12).Value = 600
$ws.Cells.Item(116, 13).Value = 550
$ws.Cells.Item(116, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(116, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(116, 16).Value = 275
$ws.Cells.Item(116, 17).Value = 2
$ws.Cells.Item(116, 18).Value = "Hortaliza"

# New row 117: "Segunda" quality sample
$ws.Cells.Item(117, 1).Value = 8
$ws.Cells.Item(117, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(117, 3).Value = "Coquimbo"
$ws.Cells.Item(117, 4).Value = 44494
$ws.Cells.Item(117, 5).Value = 4
$ws.Cells.Item(117, 6).Value = 100112009
$ws.Cells.Item(117, 7).Value = "Acelga"
$ws.Cells.Item(117, 8).Value = "Sin especificar"
$ws.Cells.Item(117, 9).Value = "Segunda"
$ws.Cells.Item(117, 10).Value = 1440
$ws.Cells.Item(117, 11).Value = 400
$ws.Cells.Item(117, 12).Value = 450
$ws.Cells.Item(117, 13).Value = 425
$ws.Cells.Item(117, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(117, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(117, 16).Value = 212
$ws.Cells.Item(117, 17).Value = 2
$ws.Cells.Item(117, 18).Value = "Hortaliza"
